$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark from the "JS is single threaded..." paragraph
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Replace the last paragraph ("fetch(url).then... ") so that the stray
#    <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> inside its <w:pPr> is
#    dropped, while the run-level formatting is preserved, and then append
#    all of the new paragraphs describing state vs props.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2E1F4C03" w14:textId="7721C0AA" w:rsidR="00601D75" w:rsidRPr="00802D6F" w:rsidRDefault="009D475B" w:rsidP="009D475B">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>f</w:t>
  </w:r>
  <w:r>
    <w:t>etch(url).then… equals to await fetch(url).</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:ind w:firstLineChars="0"/>
  </w:pPr>
  <w:r>
    <w:t>S</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>tate</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>an</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">d props: </w:t>
  </w:r>
  <w:r>
    <w:t>props and state are related. The state of one component will often become the props of a child component. Props are passed to the child within the render method of the parent as the second argument to React.createElement. State is equivalent to local variables in a function, props on the other hand, make comopnets reusable by giving components the ability to receive data from their parent component in the form of props. It is equivalent to function parameters.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:t>Class SampleComponent extends React.Component{</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:t>Render(){</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:t>Return &lt;div&gt;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> Hello {this.props.name}&lt;/div&gt;;</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:t>}</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:t>}</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>&lt;</w:t>
  </w:r>
  <w:r>
    <w:t>SampleComponent name=”Joni” /&gt;</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t>&lt;SampleComponent name =”Wang” /&gt;</w:t>
  </w:r>
</w:p>
"@

$lastRange.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) Re-add the _GoBack bookmark at the very end of the document (end of the
#    last new paragraph's text, right after the "Wang" run).
# ---------------------------------------------------------------------------
$finalPara = $d.Paragraphs.Last
$finalRange = $finalPara.Range
$finalRange.Collapse(0)
$finalRange.MoveEnd(1, -1) | Out-Null
$finalRange.MoveStart(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $finalRange) | Out-Null
